$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 7, pushing the existing rows 7-16 down to 9-18
$ws.Rows("7:8").Insert()

# New row 7: Brooks / Primera, $/bandeja 10 kilos, Región de O'Higgins
$ws.Cells.Item(7, 1).Value = 1
$ws.Cells.Item(7, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(7, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(7, 4).Value = 44914
$ws.Cells.Item(7, 5).Value = 15
$ws.Cells.Item(7, 6).Value = "Fruta"
$ws.Cells.Item(7, 7).Value = 100103
$ws.Cells.Item(7, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(7, 9).Value = 100103001
$ws.Cells.Item(7, 10).Value = "Cereza"
$ws.Cells.Item(7, 11).Value = "Brooks"
$ws.Cells.Item(7, 12).Value = "Primera"
$ws.Cells.Item(7, 13).Value = 700
$ws.Cells.Item(7, 14).Value = 7000
$ws.Cells.Item(7, 15).Value = 8000
$ws.Cells.Item(7, 16).Value = 7429
$ws.Cells.Item(7, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(7, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(7, 19).Value = 743
$ws.Cells.Item(7, 20).Value = 10

# New row 8: Lapins / Primera, $/bandeja 10 kilos, Región de O'Higgins
$ws.Cells.Item(8, 1).Value = 1
$ws.Cells.Item(8, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(8, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(8, 4).Value = 44914
$ws.Cells.Item(8, 5).Value = 15
$ws.Cells.Item(8, 6).Value = "Fruta"
$ws.Cells.Item(8, 7).Value = 100103
$ws.Cells.Item(8, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(8, 9).Value = 100103001
$ws.Cells.Item(8, 10).Value = "Cereza"
$ws.Cells.Item(8, 11).Value = "Lapins"
$ws.Cells.Item(8, 12).Value = "Primera"
$ws.Cells.Item(8, 13).Value = 550
$ws.Cells.Item(8, 14).Value = 7000
$ws.Cells.Item(8, 15).Value = 8000
$ws.Cells.Item(8, 16).Value = 7455
$ws.Cells.Item(8, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(8, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(8, 19).Value = 746
$ws.Cells.Item(8, 20).Value = 10

# Ensure the date column keeps the same date number format used by the rest of
# column D (style index 2 in the original workbook).
$ws.Cells.Item(7, 4).NumberFormat = $ws.Cells.Item(9, 4).NumberFormat
$ws.Cells.Item(8, 4).NumberFormat = $ws.Cells.Item(9, 4).NumberFormat
